$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "No problems" row (old row 4). Excel naturally shifts every row
# below it up by one, shrinks the A2:A9 merged region to A2:A8, and recomputes
# the sheet dimension - matching the rest of the diff without touching any
# cell styles.
$ws.Rows.Item(4).Delete()

# The two rows that used to report "MathML" (old C2=6, old C4=2 for "No
# problems") are consolidated into a single MathML count of 8.
$ws.Range("C2").Value = 8
